$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update season_ending_year_y (column O) values for specific rows
$ws.Range("O4").Value = 1909
$ws.Range("O6").Value = 424
$ws.Range("O7").Value = 292
$ws.Range("O8").Value = 3606
$ws.Range("O9").Value = 4644
$ws.Range("O10").Value = 3606
$ws.Range("O11").Value = 1007
